# "Close to ready for CRAN" - update the httk benchmark row for version 2.4.0
# (row 26 of Table1 on Sheet1) with the latest benchmark re-run numbers, and
# record a note about the underlying change (switch to chem props from ctxR).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B26").Value = 1021

$ws.Range("F26").Value = 0.9477
$ws.Range("G26").Value = 353
$ws.Range("H26").Value = 0.2716
$ws.Range("I26").Value = 353
$ws.Range("J26").Value = 1.508
$ws.Range("K26").Value = 36
$ws.Range("L26").Value = 0.9698
$ws.Range("M26").Value = 80
$ws.Range("N26").Value = 1.132
$ws.Range("O26").Value = 80
$ws.Range("P26").Value = 0.6466

$ws.Range("R26").Value = "Switched to chem props from ctxR"

# Restore the view state that Excel records after scrolling/selecting near
# the bottom of the table.
$ws.Activate()
$ws.Range("F27").Select()
